$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Range("B2").Value = "66"
